# Update localization status for file b36aa18e-0adf-4c3e-b324-85af3c25eda7.md
# from "Ready for handoff" to "In Translation" across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: zh-cn (E) and de-de (F) status columns for row 5
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"

# zh-cn sheet: Status column (C) for row 5
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"

# de-de sheet: Status column (C) for row 5
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
